# Update "want to go" counts (column F) across the 展览 / 演出 / 全部类型 sheets
# to reflect a fresh scrape of attendance numbers (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # exhibition sheet
$ws2 = $wb.Worksheets.Item(2)   # performance sheet
$ws4 = $wb.Worksheets.Item(4)   # all-types sheet

# exhibition sheet
$ws1.Range("F2").Value = 88
$ws1.Range("F3").Value = 176
$ws1.Range("F5").Value = 59
$ws1.Range("F6").Value = 2766
$ws1.Range("F7").Value = 1652
$ws1.Range("F8").Value = 7491
$ws1.Range("F10").Value = 7677
$ws1.Range("F11").Value = 22
$ws1.Range("F12").Value = 37
$ws1.Range("F13").Value = 13
$ws1.Range("F14").Value = 6240
$ws1.Range("F15").Value = 3285
$ws1.Range("F17").Value = 3648
$ws1.Range("F18").Value = 24
$ws1.Range("F19").Value = 19
$ws1.Range("F20").Value = 20
$ws1.Range("F21").Value = 37
$ws1.Range("F22").Value = 2
$ws1.Range("F26").Value = 293
$ws1.Range("F27").Value = 3673
$ws1.Range("F32").Value = 1317
$ws1.Range("F33").Value = 65
$ws1.Range("F34").Value = 22
$ws1.Range("F35").Value = 2638
$ws1.Range("F36").Value = 1537
$ws1.Range("F38").Value = 25
$ws1.Range("F39").Value = 31
$ws1.Range("F40").Value = 3336
$ws1.Range("F41").Value = 190
$ws1.Range("F42").Value = 257
$ws1.Range("F43").Value = 33
$ws1.Range("F45").Value = 489
$ws1.Range("F46").Value = 1305
$ws1.Range("F48").Value = 529
$ws1.Range("F49").Value = 601

# performance sheet
$ws2.Range("F7").Value = 21
$ws2.Range("F9").Value = 402
$ws2.Range("F13").Value = 19

# all-types sheet
$ws4.Range("F3").Value = 88
$ws4.Range("F4").Value = 176
$ws4.Range("F6").Value = 59
$ws4.Range("F7").Value = 124
$ws4.Range("F8").Value = 2766
$ws4.Range("F9").Value = 1652
$ws4.Range("F13").Value = 7491
$ws4.Range("F14").Value = 7677
$ws4.Range("F15").Value = 22
$ws4.Range("F16").Value = 37
$ws4.Range("F17").Value = 6240
$ws4.Range("F18").Value = 3285
$ws4.Range("F19").Value = 3648
$ws4.Range("F20").Value = 19
$ws4.Range("F21").Value = 37
$ws4.Range("F25").Value = 293
$ws4.Range("F26").Value = 3673
$ws4.Range("F32").Value = 1317
$ws4.Range("F33").Value = 65
$ws4.Range("F34").Value = 22
$ws4.Range("F35").Value = 2638
$ws4.Range("F36").Value = 1537
$ws4.Range("F38").Value = 25
$ws4.Range("F39").Value = 31
$ws4.Range("F40").Value = 3336
$ws4.Range("F41").Value = 190
$ws4.Range("F42").Value = 257
$ws4.Range("F43").Value = 33
$ws4.Range("F45").Value = 489
$ws4.Range("F46").Value = 1305
$ws4.Range("F48").Value = 529
$ws4.Range("F49").Value = 601

Write-Host "Done updating want-to-go counts"
